$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 773.8333
$ws.Range("I11").Value = 773.8333
$ws.Range("K11").Value = 773.8333
$ws.Range("M11").Value = -633.8333
$ws.Range("H62").Value = 38470868
$ws.Range("I62").Value = 1810.5555
$ws.Range("J62").Value = 125026250
$ws.Range("K62").Value = 1810.5555
$ws.Range("L62").Value = 125026250
$ws.Range("M62").Value = -1186.5555
$ws.Range("N62").Value = -125027498
$ws.Range("H65").Value = 38470868
$ws.Range("I65").Value = 1810.5555
$ws.Range("J65").Value = 125026250
$ws.Range("K65").Value = 9052.7775
$ws.Range("L65").Value = 625131250
$ws.Range("M65").Value = -5932.7775
$ws.Range("N65").Value = -625137490
$ws.Range("H103").Value = 8000661
$ws.Range("I103").Value = 516.0714
$ws.Range("J103").Value = 18182664
$ws.Range("K103").Value = 1548.2142
$ws.Range("L103").Value = 54547992
$ws.Range("M103").Value = -962.2142000000001
$ws.Range("N103").Value = -54549164
$ws.Range("H132").Value = 2197.1528
$ws.Range("I132").Value = 1492.5454
$ws.Range("K132").Value = 4477.6362
$ws.Range("M132").Value = -1947.6362
$ws.Range("H136").Value = 83746.36
$ws.Range("J136").Value = 83746.36
$ws.Range("L136").Value = 83746.36
$ws.Range("N136").Value = -93946.36

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23950.408
$ws.Range("I32").Value = 8901.154
$ws.Range("J32").Value = 102206.53
$ws.Range("K32").Value = 8901.154
$ws.Range("L32").Value = 102206.53
$ws.Range("M32").Value = -8614.154
$ws.Range("N32").Value = -102780.53
$ws.Range("H45").Value = 4067.3333
$ws.Range("I45").Value = 2884.3333
$ws.Range("J45").Value = 6433.3335
$ws.Range("K45").Value = 2884.3333
$ws.Range("L45").Value = 6433.3335
$ws.Range("M45").Value = -2507.3333
$ws.Range("N45").Value = -7187.3335
$ws.Range("H74").Value = 20680.81
$ws.Range("I74").Value = 1161.8889
$ws.Range("J74").Value = 130474.75
$ws.Range("K74").Value = 1161.8889
$ws.Range("L74").Value = 130474.75
$ws.Range("M74").Value = -287.8888999999999
$ws.Range("N74").Value = -132222.75
$ws.Range("H77").Value = 20680.81
$ws.Range("I77").Value = 1161.8889
$ws.Range("J77").Value = 130474.75
$ws.Range("K77").Value = 5809.4445
$ws.Range("L77").Value = 652373.75
$ws.Range("M77").Value = -1441.4445
$ws.Range("N77").Value = -661109.75
$ws.Range("H122").Value = 27874
$ws.Range("I122").Value = 51128
$ws.Range("K122").Value = 153384
$ws.Range("M122").Value = -150934
$ws.Range("H132").Value = 2172.8857
$ws.Range("I132").Value = 1868.0952
$ws.Range("J132").Value = 2630.0715
$ws.Range("K132").Value = 5604.2856
$ws.Range("L132").Value = 7890.2145
$ws.Range("M132").Value = -3074.2856
$ws.Range("N132").Value = -12950.2145

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 33106.855
$ws.Range("J35").Value = 33106.855
$ws.Range("L35").Value = 33106.855
$ws.Range("N35").Value = -33726.855
$ws.Range("H82").Value = 23818.65
$ws.Range("I82").Value = 8476.333000000001
$ws.Range("J82").Value = 36371.453
$ws.Range("K82").Value = 8476.333000000001
$ws.Range("L82").Value = 36371.453
$ws.Range("M82").Value = -8093.333000000001
$ws.Range("N82").Value = -37137.453
$ws.Range("H85").Value = 23818.65
$ws.Range("I85").Value = 8476.333000000001
$ws.Range("J85").Value = 36371.453
$ws.Range("K85").Value = 8476.333000000001
$ws.Range("L85").Value = 36371.453
$ws.Range("M85").Value = -7150.333000000001
$ws.Range("N85").Value = -39023.453

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 17259.334
$ws.Range("J68").Value = 17259.334
$ws.Range("L68").Value = 17259.334
$ws.Range("N68").Value = -18757.334
$ws.Range("H71").Value = 17259.334
$ws.Range("J71").Value = 17259.334
$ws.Range("L71").Value = 51778.00199999999
$ws.Range("N71").Value = -59266.00199999999
$ws.Range("H122").Value = 1751.625
$ws.Range("I122").Value = 1401
$ws.Range("J122").Value = 2803.5
$ws.Range("K122").Value = 4203
$ws.Range("L122").Value = 8410.5
$ws.Range("M122").Value = -1753
$ws.Range("N122").Value = -13310.5
$ws.Range("H134").Value = 3640.182
$ws.Range("I134").Value = 4027.7368
$ws.Range("J134").Value = 1185.6666
$ws.Range("K134").Value = 12083.2104
$ws.Range("L134").Value = 3556.9998
$ws.Range("M134").Value = -9548.2104
$ws.Range("N134").Value = -8626.9998

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1033
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1033
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3099
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4721
$ws.Range("H71").Value = 1033
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1033
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 9297
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -17409
$ws.Range("H131").Value = 915.2449
$ws.Range("J131").Value = 968.6818
$ws.Range("L131").Value = 2906.0454
$ws.Range("N131").Value = -12986.0454
$ws.Range("H132").Value = 289621.56
$ws.Range("I132").Value = 775141.75
$ws.Range("J132").Value = 5006.276
$ws.Range("K132").Value = 6976275.75
$ws.Range("L132").Value = 45056.484
$ws.Range("M132").Value = -6973745.75
$ws.Range("N132").Value = -50116.484

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 30658.719
$ws.Range("I70").Value = 34510.94
$ws.Range("K70").Value = 34510.94
$ws.Range("M70").Value = -34240.94
$ws.Range("H73").Value = 30658.719
$ws.Range("I73").Value = 34510.94
$ws.Range("K73").Value = 34510.94
$ws.Range("M73").Value = -33574.94
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1282.8572
$ws.Range("I22").Value = 826.6667
$ws.Range("J22").Value = 1625
$ws.Range("K22").Value = 826.6667
$ws.Range("L22").Value = 1625
$ws.Range("M22").Value = -531.6667
$ws.Range("N22").Value = -2215
$ws.Range("H27").Value = 1282.8572
$ws.Range("I27").Value = 826.6667
$ws.Range("J27").Value = 1625
$ws.Range("K27").Value = 826.6667
$ws.Range("L27").Value = 1625
$ws.Range("M27").Value = -719.6667
$ws.Range("N27").Value = -1839
$ws.Range("H122").Value = 4235.706
$ws.Range("I122").Value = 3800.75
$ws.Range("J122").Value = 4369.5386
$ws.Range("K122").Value = 11402.25
$ws.Range("L122").Value = 13108.6158
$ws.Range("M122").Value = -8952.25
$ws.Range("N122").Value = -18008.6158
$ws.Range("H136").Value = 3735.8445
$ws.Range("I136").Value = 1377.7222
$ws.Range("J136").Value = 13168.333
$ws.Range("K136").Value = 4133.1666
$ws.Range("L136").Value = 39504.999
$ws.Range("M136").Value = -1583.1666
$ws.Range("N136").Value = -44604.999

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 9900
$ws.Range("J47").Value = 9900
$ws.Range("L47").Value = 9900
$ws.Range("N47").Value = -11044
$ws.Range("H109").Value = 13938.5
$ws.Range("J109").Value = 13938.5
$ws.Range("L109").Value = 13938.5
$ws.Range("N109").Value = -16712.5
$ws.Range("H122").Value = 7306.5127
$ws.Range("I122").Value = 9597.360000000001
$ws.Range("J122").Value = 3215.7144
$ws.Range("K122").Value = 28792.08
$ws.Range("L122").Value = 9647.143199999999
$ws.Range("M122").Value = -26342.08
$ws.Range("N122").Value = -14547.1432
$ws.Range("H132").Value = 3908.04
$ws.Range("I132").Value = 4188.706
$ws.Range("J132").Value = 3311.625
$ws.Range("K132").Value = 12566.118
$ws.Range("L132").Value = 9934.875
$ws.Range("M132").Value = -10036.118
$ws.Range("N132").Value = -14994.875
